$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Content updates on "List of Predictors" (sheet 1) ---
# Row 8: the "Not sure..." note gets an extra sentence about Corine 1990s data.
$ws1.Range("G8").Value = "Not sure if the data is available in a usable format.. Could maybe be combined with Land use change data for CZ? Does it indicate the type of habitat? Or remote sensing land cover data? Else: Corine 1990s data (collected 85-86)"

# Row 9: hypothesis reworded, reference swapped for a DOI link, CHELSA note gets a double space typo.
$ws1.Range("D9").Value = "Climate (+climate geography) can predict temporal change"
$ws1.Range("E9").Value = "https://doi.org/10.1038/s41586-023-06577-5"
$ws1.Range("F9").Value = "CHELSA?  - for 80s or 2000s time period"

# --- View / selection updates ---
# Sheet1: drop the frozen/scrolled topLeftCell, select E4 instead of E23, keep sheet1 active.
$ws1.Activate()
$ws1.Range("E4").Select() | Out-Null

# Sheet2: change selection to G18.
$ws2.Activate()
$ws2.Range("G18").Select() | Out-Null

# Re-activate sheet1 so it remains the selected tab on save.
$ws1.Activate()

# --- Row height follow-up for the longer G8 text (auto-fit style growth) ---
$ws1.Rows.Item(8).RowHeight = 105
